$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 100 -> 0M
$t.Cell(1,1).Range.Text = "0M"

# Row 2: 0 -> 0M
$t.Cell(2,1).Range.Text = "0M"

# Row 3: 228 -> 0M
$t.Cell(3,1).Range.Text = "0M"

# Row 4: 3 -> 55
$t.Cell(4,1).Range.Text = "55"

# Row 5: 0.00004 -> 0.00003
$t.Cell(5,1).Range.Text = "0.00003"

# Insert a new row after row 5 (before the current row 6) with text 0.00012
$beforeRow = $t.Rows.Item(6)
$newRow = $t.Rows.Add($beforeRow)
$t.Cell(6,1).Range.Text = "0.00012"

# Remove the row that now holds 0.00005 (was between 0.00006 and 0.00001)
$t.Rows.Item(8).Delete()

# Row (now 9): 0.00004 -> 0.00008
$t.Cell(9,1).Range.Text = "0.00008"

# Row (now 10): 0.00005 -> 0.00008
$t.Cell(10,1).Range.Text = "0.00008"

# Row (now 11): 0.00006 -> 0.00010
$t.Cell(11,1).Range.Text = "0.00010"

# Row (now 12): 0.00016 -> 0.00444
$t.Cell(12,1).Range.Text = "0.00444"

# Row 44 (tab-separated run starting with "2") -> collapse to single value "100"
$t.Cell(44,1).Range.Text = "100"

# Row 45 (tab-separated run starting with "40") -> collapse to single value "0"
$t.Cell(45,1).Range.Text = "0"

# Row 46 (tab-separated run starting with "10") -> collapse to single value "228"
$t.Cell(46,1).Range.Text = "228"
